# Update team-specific transition matrix values (Sheet1) to reflect
# newly added/updated team-specific time data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "B2"; Value = 0.1780821917808219 },
    @{ Cell = "C2"; Value = 0.6027397260273972 },
    @{ Cell = "J2"; Value = 0.02465753424657534 },
    @{ Cell = "P2"; Value = 0.1232876712328767 },
    @{ Cell = "S2"; Value = 0.07123287671232877 },
    @{ Cell = "C3"; Value = 0.0045662100456621 },
    @{ Cell = "J3"; Value = 0.0730593607305936 },
    @{ Cell = "P3"; Value = 0.7397260273972602 },
    @{ Cell = "S3"; Value = 0.182648401826484 },
    @{ Cell = "J4"; Value = 0.1052631578947368 },
    @{ Cell = "P4"; Value = 0.7105263157894737 },
    @{ Cell = "S4"; Value = 0.1842105263157895 },
    @{ Cell = "B6"; Value = 0.05384615384615385 },
    @{ Cell = "D6"; Value = 0.007692307692307693 },
    @{ Cell = "F6"; Value = 0.06153846153846154 },
    @{ Cell = "J6"; Value = 0.2807692307692308 },
    @{ Cell = "O6"; Value = 0.01923076923076923 },
    @{ Cell = "Q6"; Value = 0.1769230769230769 },
    @{ Cell = "R6"; Value = 0.08076923076923077 },
    @{ Cell = "S6"; Value = 0.3192307692307692 },
    @{ Cell = "B7"; Value = 0.09615384615384616 },
    @{ Cell = "D7"; Value = 0.009615384615384616 },
    @{ Cell = "E7"; Value = 0.004807692307692308 },
    @{ Cell = "F7"; Value = 0.04807692307692308 },
    @{ Cell = "J7"; Value = 0.125 },
    @{ Cell = "O7"; Value = 0.01923076923076923 },
    @{ Cell = "Q7"; Value = 0.1778846153846154 },
    @{ Cell = "R7"; Value = 0.1201923076923077 },
    @{ Cell = "S7"; Value = 0.3990384615384616 },
    @{ Cell = "B8"; Value = 0.09881422924901186 },
    @{ Cell = "D8"; Value = 0.0158102766798419 },
    @{ Cell = "E8"; Value = 0.001976284584980237 },
    @{ Cell = "F8"; Value = 0.06719367588932806 },
    @{ Cell = "J8"; Value = 0.1126482213438735 },
    @{ Cell = "O8"; Value = 0.01976284584980237 },
    @{ Cell = "Q8"; Value = 0.1600790513833992 },
    @{ Cell = "R8"; Value = 0.1067193675889328 },
    @{ Cell = "S8"; Value = 0.41699604743083 },
    @{ Cell = "B9"; Value = 0.07468879668049792 },
    @{ Cell = "D9"; Value = 0.008298755186721992 },
    @{ Cell = "F9"; Value = 0.04979253112033195 },
    @{ Cell = "J9"; Value = 0.09958506224066389 },
    @{ Cell = "O9"; Value = 0.01244813278008299 },
    @{ Cell = "Q9"; Value = 0.1701244813278008 },
    @{ Cell = "R9"; Value = 0.1369294605809129 },
    @{ Cell = "S9"; Value = 0.4481327800829876 },
    @{ Cell = "B10"; Value = 0.1214421252371917 },
    @{ Cell = "D10"; Value = 0.01644528779253637 },
    @{ Cell = "E10"; Value = 0.0006325110689437065 },
    @{ Cell = "F10"; Value = 0.06451612903225806 },
    @{ Cell = "J10"; Value = 0.1347248576850095 },
    @{ Cell = "O10"; Value = 0.01644528779253637 },
    @{ Cell = "Q10"; Value = 0.200506008855155 },
    @{ Cell = "R10"; Value = 0.1005692599620493 },
    @{ Cell = "S10"; Value = 0.3447185325743201 },
    @{ Cell = "G11"; Value = 0.1437699680511182 },
    @{ Cell = "J11"; Value = 0.0670926517571885 },
    @{ Cell = "K11"; Value = 0.1980830670926517 },
    @{ Cell = "L11"; Value = 0.5782747603833865 },
    @{ Cell = "S11"; Value = 0.01277955271565495 },
    @{ Cell = "G12"; Value = 0.7128205128205128 },
    @{ Cell = "J12"; Value = 0.2102564102564103 },
    @{ Cell = "K12"; Value = 0.01538461538461539 },
    @{ Cell = "L12"; Value = 0.03589743589743589 },
    @{ Cell = "S12"; Value = 0.02564102564102564 },
    @{ Cell = "G13"; Value = 0.7441860465116279 },
    @{ Cell = "J13"; Value = 0.2325581395348837 },
    @{ Cell = "S13"; Value = 0.02325581395348837 },
    @{ Cell = "G14"; Value = 0.6666666666666666 },
    @{ Cell = "J14"; Value = 0.3333333333333333 },
    @{ Cell = "H15"; Value = 0.144578313253012 },
    @{ Cell = "I15"; Value = 0.05622489959839357 },
    @{ Cell = "J15"; Value = 0.4257028112449799 },
    @{ Cell = "K15"; Value = 0.04417670682730924 },
    @{ Cell = "M15"; Value = 0.008032128514056224 },
    @{ Cell = "N15"; Value = 0.004016064257028112 },
    @{ Cell = "O15"; Value = 0.07630522088353414 },
    @{ Cell = "S15"; Value = 0.2409638554216867 },
    @{ Cell = "F16"; Value = 0.01754385964912281 },
    @{ Cell = "H16"; Value = 0.1052631578947368 },
    @{ Cell = "I16"; Value = 0.07456140350877193 },
    @{ Cell = "J16"; Value = 0.5043859649122807 },
    @{ Cell = "K16"; Value = 0.06140350877192982 },
    @{ Cell = "M16"; Value = 0.02192982456140351 },
    @{ Cell = "N16"; Value = 0.004385964912280702 },
    @{ Cell = "O16"; Value = 0.07456140350877193 },
    @{ Cell = "S16"; Value = 0.1359649122807018 },
    @{ Cell = "F17"; Value = 0.02912621359223301 },
    @{ Cell = "H17"; Value = 0.2097087378640777 },
    @{ Cell = "I17"; Value = 0.09320388349514563 },
    @{ Cell = "J17"; Value = 0.3689320388349515 },
    @{ Cell = "K17"; Value = 0.09320388349514563 },
    @{ Cell = "M17"; Value = 0.01941747572815534 },
    @{ Cell = "N17"; Value = 0.001941747572815534 },
    @{ Cell = "O17"; Value = 0.06407766990291262 },
    @{ Cell = "S17"; Value = 0.1203883495145631 },
    @{ Cell = "F18"; Value = 0.02422145328719723 },
    @{ Cell = "H18"; Value = 0.1660899653979239 },
    @{ Cell = "I18"; Value = 0.08996539792387544 },
    @{ Cell = "J18"; Value = 0.4567474048442907 },
    @{ Cell = "K18"; Value = 0.1038062283737024 },
    @{ Cell = "M18"; Value = 0.006920415224913495 },
    @{ Cell = "N18"; Value = 0.01038062283737024 },
    @{ Cell = "O18"; Value = 0.03114186851211072 },
    @{ Cell = "S18"; Value = 0.1107266435986159 },
    @{ Cell = "F19"; Value = 0.01642710472279261 },
    @{ Cell = "H19"; Value = 0.1923340177960301 },
    @{ Cell = "I19"; Value = 0.09308692676249145 },
    @{ Cell = "J19"; Value = 0.3894592744695414 },
    @{ Cell = "K19"; Value = 0.09582477754962354 },
    @{ Cell = "M19"; Value = 0.01779603011635866 },
    @{ Cell = "O19"; Value = 0.06639288158795345 },
    @{ Cell = "S19"; Value = 0.1286789869952088 }
)

foreach ($update in $updates) {
    $ws.Range($update.Cell).Value = $update.Value
}
